$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: E7 goes from 0 to 0.75 hours, formatted like D7 (red, centered text)
$ws.Range("E7").Value = 0.75
$ws.Range("E7").Font.Color = $ws.Range("D7").Font.Color()
$ws.Range("E7").HorizontalAlignment = $ws.Range("D7").HorizontalAlignment()
$ws.Range("E7").VerticalAlignment = $ws.Range("D7").VerticalAlignment()

# Row 10: E10 goes from 0 to 0.25 hours, same red/centered formatting
$ws.Range("E10").Value = 0.25
$ws.Range("E10").Font.Color = $ws.Range("D7").Font.Color()
$ws.Range("E10").HorizontalAlignment = $ws.Range("D7").HorizontalAlignment()
$ws.Range("E10").VerticalAlignment = $ws.Range("D7").VerticalAlignment()

# Move the active cell selection from D7 to E7
$ws.Range("E7").Select()
